$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column B (Expected Time) values 30 -> 47
$ws.Range("B1").Value = 47
$ws.Range("B2").Value = 47
$ws.Range("B3").Value = 47
$ws.Range("B4").Value = 47

# Clear column D rows 2-4 (Team Progress actuals), keep D1 = 0
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()

# Update the active selection to D8 as recorded when file was last saved
$ws.Range("D8").Select()
